$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row -> new "F" value (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7021
$ws1.Range("F4").Value = 461
$ws1.Range("F7").Value = 147
$ws1.Range("F13").Value = 445
$ws1.Range("F15").Value = 1827
$ws1.Range("F16").Value = 44
$ws1.Range("F17").Value = 3626
$ws1.Range("F18").Value = 26
$ws1.Range("F21").Value = 23
$ws1.Range("F22").Value = 27
$ws1.Range("F23").Value = 2254
$ws1.Range("F24").Value = 16
$ws1.Range("F25").Value = 252
$ws1.Range("F30").Value = 18
$ws1.Range("F32").Value = 240
$ws1.Range("F33").Value = 92

# Sheet "全部类型" (All types) - row -> new "F" value (想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7021
$ws4.Range("F4").Value = 461
$ws4.Range("F8").Value = 147
$ws4.Range("F14").Value = 445
$ws4.Range("F16").Value = 1827
$ws4.Range("F17").Value = 44
$ws4.Range("F18").Value = 3626
$ws4.Range("F19").Value = 26
$ws4.Range("F22").Value = 23
$ws4.Range("F23").Value = 27
$ws4.Range("F24").Value = 2254
$ws4.Range("F25").Value = 16
$ws4.Range("F26").Value = 252
$ws4.Range("F31").Value = 18
$ws4.Range("F33").Value = 240
$ws4.Range("F34").Value = 92
